# Auto-generated edit script applying the crypto price/volume refresh
# described in the commit "Updated cryptos list on Tue Aug 27 09:36:02 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.893.76"
$ws.Range("E2").Value = "  -1.14%  "
$ws.Range("D3").Value = "2.681.95"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "555.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.589"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.12%  "
$ws.Range("E9").Value = "  -3.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.162"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.50"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.98%  "
$ws.Range("E12").Value = "  -3.14%  "
$ws.Range("D13").Value = "3.157.30"
$ws.Range("E13").Value = "  -1.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.77%  "
$ws.Range("D15").Value = "62.802.36"
$ws.Range("E15").Value = "  -0.97%  "
$ws.Range("E16").Value = "  -1.94%  "
$ws.Range("D17").Value = "2.683.16"
$ws.Range("E17").Value = "  -1.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "345.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.26%  "
$ws.Range("E21").Value = "  -5.04%  "
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.509"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.44%  "
$ws.Range("E25").Value = "  +0.26%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.11%  "
$ws.Range("E28").Value = "  +5.58%  "
$ws.Range("D29").Value = "0.0₃0852"
$ws.Range("E29").Value = "  -5.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.99%  "
$ws.Range("E31").Value = "  -1.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "163.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.69%  "
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("E35").Value = "  -0.54%  "
$ws.Range("E36").Value = "  -2.67%  "
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "340.65"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.98%  "
$ws.Range("E39").Value = "  -1.46%  "
$ws.Range("E40").Value = "  -4.27%  "
$ws.Range("E41").Value = "  -1.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.81%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.615"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.46%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0557"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.15%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.84"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.37%  "
$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0970"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.69%  "
$ws.Range("E51").Value = "  -3.28%  "
